$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 12422.543
$ws.Range("I15").Value = 12422.543
$ws.Range("K15").Value = 37267.629
$ws.Range("M15").Value = -37098.629

$ws.Range("H98").Value = 2123.182
$ws.Range("I98").Value = 2150.5557
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 2150.5557
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = -652.5556999999999
$ws.Range("N98").Value = -4996

$ws.Range("H122").Value = 2123.182
$ws.Range("I122").Value = 2150.5557
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 6451.6671
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -4001.6671
$ws.Range("N122").Value = -10900

$ws.Range("H137").Value = 1288.9482
$ws.Range("I137").Value = 730.0833
$ws.Range("J137").Value = 1683.4412
$ws.Range("K137").Value = 2190.2499
$ws.Range("L137").Value = 5050.3236
$ws.Range("M137").Value = 359.7501000000002
$ws.Range("N137").Value = -10150.3236

$ws.Range("H138").Value = 3633.6448
$ws.Range("I138").Value = 1459.1072
$ws.Range("J138").Value = 4902.125
$ws.Range("K138").Value = 4377.321599999999
$ws.Range("L138").Value = 14706.375
$ws.Range("M138").Value = 762.6784000000007
$ws.Range("N138").Value = -24986.375

$ws.Range("H141").Value = 2032.96
$ws.Range("I141").Value = 2032.96
$ws.Range("K141").Value = 6098.88
$ws.Range("M141").Value = -918.8800000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 250001950
$ws.Range("I88").Value = 3000
$ws.Range("J88").Value = 333334940
$ws.Range("K88").Value = 3000
$ws.Range("L88").Value = 333334940
$ws.Range("M88").Value = -2594
$ws.Range("N88").Value = -333335752

$ws.Range("H91").Value = 250001950
$ws.Range("I91").Value = 3000
$ws.Range("J91").Value = 333334940
$ws.Range("K91").Value = 3000
$ws.Range("L91").Value = 333334940
$ws.Range("M91").Value = -1596
$ws.Range("N91").Value = -333337748

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 28828.46
$ws.Range("J126").Value = 28828.46
$ws.Range("L126").Value = 28828.46
$ws.Range("N126").Value = -38708.46

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2526.5454
$ws.Range("I31").Value = 2050.85
$ws.Range("J31").Value = 2666.4558
$ws.Range("K31").Value = 2050.85
$ws.Range("L31").Value = 2666.4558
$ws.Range("M31").Value = -1755.85
$ws.Range("N31").Value = -3256.4558

$ws.Range("H33").Value = 8020.3335
$ws.Range("I33").Value = 4530.5
$ws.Range("J33").Value = 15000
$ws.Range("K33").Value = 4530.5
$ws.Range("L33").Value = 15000
$ws.Range("M33").Value = -4151.5
$ws.Range("N33").Value = -15758

$ws.Range("H34").Value = 2526.5454
$ws.Range("I34").Value = 2050.85
$ws.Range("J34").Value = 2666.4558
$ws.Range("K34").Value = 2050.85
$ws.Range("L34").Value = 2666.4558
$ws.Range("M34").Value = -1848.85
$ws.Range("N34").Value = -3070.4558

$ws.Range("H58").Value = 1375.6428
$ws.Range("I58").Value = 966.2353000000001
$ws.Range("J58").Value = 2008.3636
$ws.Range("K58").Value = 966.2353000000001
$ws.Range("L58").Value = 2008.3636
$ws.Range("M58").Value = -763.2353000000001
$ws.Range("N58").Value = -2414.3636

$ws.Range("H120").Value = 45800
$ws.Range("J120").Value = 45800
$ws.Range("L120").Value = 45800
$ws.Range("N120").Value = -53058

$ws.Range("H132").Value = 2142.25
$ws.Range("I132").Value = 1437.8695
$ws.Range("J132").Value = 5382.4
$ws.Range("K132").Value = 4313.6085
$ws.Range("L132").Value = 16147.2
$ws.Range("M132").Value = -1783.6085
$ws.Range("N132").Value = -21207.2

$ws.Range("H134").Value = 1413.037
$ws.Range("I134").Value = 1448.1538
$ws.Range("J134").Value = 500
$ws.Range("K134").Value = 4344.4614
$ws.Range("L134").Value = 1500
$ws.Range("M134").Value = -1809.4614
$ws.Range("N134").Value = -6570

$ws.Range("H136").Value = 1375.6428
$ws.Range("I136").Value = 966.2353000000001
$ws.Range("J136").Value = 2008.3636
$ws.Range("K136").Value = 2898.7059
$ws.Range("L136").Value = 6025.0908
$ws.Range("M136").Value = -348.7058999999999
$ws.Range("N136").Value = -11125.0908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5928.5713
$ws.Range("I56").Value = 5928.5713
$ws.Range("K56").Value = 5928.5713
$ws.Range("M56").Value = -5398.5713

$ws.Range("H68").Value = 1040.43
$ws.Range("I68").Value = 789.37036
$ws.Range("J68").Value = 1335.1522
$ws.Range("K68").Value = 2368.11108
$ws.Range("L68").Value = 4005.4566
$ws.Range("M68").Value = -1557.11108
$ws.Range("N68").Value = -5627.4566

$ws.Range("H71").Value = 1040.43
$ws.Range("I71").Value = 789.37036
$ws.Range("J71").Value = 1335.1522
$ws.Range("K71").Value = 7104.33324
$ws.Range("L71").Value = 12016.3698
$ws.Range("M71").Value = -3048.33324
$ws.Range("N71").Value = -20128.3698

$ws.Range("H113").Value = 1017444.4
$ws.Range("I113").Value = 1389342.4
$ws.Range("J113").Value = 435343.22
$ws.Range("K113").Value = 4168027.2
$ws.Range("L113").Value = 1306029.66
$ws.Range("M113").Value = -4165857.2
$ws.Range("N113").Value = -1310369.66

$ws.Range("H122").Value = 509.72223
$ws.Range("I122").Value = 448.83334
$ws.Range("J122").Value = 527.119
$ws.Range("K122").Value = 4039.50006
$ws.Range("L122").Value = 4744.071
$ws.Range("M122").Value = -1589.50006
$ws.Range("N122").Value = -9644.071

$ws.Range("H123").Value = 3000
$ws.Range("I123").Value = 3000
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 9000
$ws.Range("L123").Value = 0
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -6550

$ws.Range("H131").Value = 13102462
$ws.Range("I131").Value = 6700485.5
$ws.Range("J131").Value = 14494196
$ws.Range("K131").Value = 20101456.5
$ws.Range("L131").Value = 43482588
$ws.Range("M131").Value = -20096416.5
$ws.Range("N131").Value = -43492668

$ws.Range("H132").Value = 2370.274
$ws.Range("I132").Value = 1471.5555
$ws.Range("J132").Value = 2478.12
$ws.Range("K132").Value = 13243.9995
$ws.Range("L132").Value = 22303.08
$ws.Range("M132").Value = -10713.9995
$ws.Range("N132").Value = -27363.08

$ws.Range("H133").Value = 53475.76
$ws.Range("I133").Value = 102399.1
$ws.Range("J133").Value = 9000
$ws.Range("K133").Value = 307197.3
$ws.Range("L133").Value = 27000
$ws.Range("M133").Value = -302137.3
$ws.Range("N133").Value = -37120

$ws.Range("H134").Value = 11200.272
$ws.Range("I134").Value = 14951.125
$ws.Range("K134").Value = 44853.375
$ws.Range("M134").Value = -39783.375

$ws.Range("H137").Value = 38480700
$ws.Range("I137").Value = 3399.75
$ws.Range("J137").Value = 45476572
$ws.Range("K137").Value = 10199.25
$ws.Range("L137").Value = 136429716
$ws.Range("M137").Value = -5099.25
$ws.Range("N137").Value = -136439916

$ws.Range("H138").Value = 10860
$ws.Range("I138").Value = 10705.454
$ws.Range("J138").Value = 11200
$ws.Range("K138").Value = 32116.362
$ws.Range("L138").Value = 33600
$ws.Range("M138").Value = -26976.362
$ws.Range("N138").Value = -43880

$ws.Range("H139").Value = 5159.0312
$ws.Range("I139").Value = 7775.933
$ws.Range("J139").Value = 2850
$ws.Range("K139").Value = 23327.799
$ws.Range("L139").Value = 8550
$ws.Range("M139").Value = -18187.799
$ws.Range("N139").Value = -18830

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4390
$ws.Range("I80").Value = 4002.5
$ws.Range("J80").Value = 4777.5
$ws.Range("K80").Value = 4002.5
$ws.Range("L80").Value = 4777.5
$ws.Range("M80").Value = -3004.5
$ws.Range("N80").Value = -6773.5

$ws.Range("H83").Value = 4390
$ws.Range("I83").Value = 4002.5
$ws.Range("J83").Value = 4777.5
$ws.Range("K83").Value = 20012.5
$ws.Range("L83").Value = 23887.5
$ws.Range("M83").Value = -15020.5
$ws.Range("N83").Value = -33871.5

$ws.Range("H97").Value = 756.55554
$ws.Range("I97").Value = 756.55554
$ws.Range("K97").Value = 756.55554
$ws.Range("M97").Value = -260.55554

$ws.Range("H119").Value = 20000
$ws.Range("J119").Value = 20000
$ws.Range("L119").Value = 20000
$ws.Range("N119").Value = -29676

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H101").Value = 10874
$ws.Range("J101").Value = 10874
$ws.Range("L101").Value = 10874
$ws.Range("N101").Value = -17364

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2122.3
$ws.Range("J96").Value = 2477.7273
$ws.Range("L96").Value = 2477.7273
$ws.Range("N96").Value = -5223.7273
